# Applies the "_Release Schedule.xlsx" update:
#  - Statistics block: split "Hypothesis Testing IV/V" rows into
#    "Estimation Statistics I/II" with corrected dates.
#  - Forecasting block: date corrections.
#  - Machine Learning block: "Regression"/"Classification" rows each split
#    into two (Linear/Penalized + Nonlinear variants), existing rows shift
#    down to make room.
#  - Selection moves to E21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal-text value into a cell without Excel's
# autodetection turning date-looking strings into real dates (which would
# also stamp a date NumberFormat on the cell). We briefly mark the cell as
# Text, assign, then restore the "Normal" style so no stray formatting is
# left behind.
function Set-TextValue {
    param($addr, $value)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# ---- Statistics block -----------------------------------------------
# Probability Distributions row: date correction 10/30 -> 10/23
Set-TextValue "C5" "10/23/2022"

# Rows 9 & 10: renamed topics + date corrections
$ws.Range("B9").Value = "Estimation Statistics I - Statistical Power"
Set-TextValue "C9" "10/23/2022"

$ws.Range("B10").Value = "Estimation Statistics II - Effect Size"
Set-TextValue "C10" "10/23/2022"

# ---- Forecasting block ------------------------------------------------
Set-TextValue "C12" "10/30/2022"
Set-TextValue "C13" "11/06/2022"
Set-TextValue "C14" "11/13/2022"

# ---- Machine Learning block -------------------------------------------
# "Regression" (row 16) and "Classification" (row 17) each split into two
# rows. Insert a blank row right after each so the rows below (Decision
# Trees, Ensembles, Dimension Reduction, Unsupervised Learning, Deep
# Learning) shift down intact, then fill in the four new topic cells.
$ws.Rows(17).Insert() | Out-Null
$ws.Rows(19).Insert() | Out-Null

$ws.Range("B16").Value = "Regression I - Linear and Penalized Regression"

$ws.Range("A17").Value = "Machine Learning"
$ws.Range("B17").Value = "Regression II - Nonlinear Regression Models"
$ws.Range("C17").Value = "TBD"

$ws.Range("B18").Value = "Classification I - Linear and Penalized Classification"

$ws.Range("A19").Value = "Machine Learning"
$ws.Range("B19").Value = "Classification II - Nonlinear Classification Models"
$ws.Range("C19").Value = "TBD"

# ---- Selection ----------------------------------------------------------
$ws.Range("E21").Select() | Out-Null
